# "finestra incidenza 7gg centrata su ultimo g"
#
# The sheet tracks daily new positives (column B, "nuovi pos.") and a
# 7-day rolling sum (column C, "somma mobile 7gg.") together with its
# per-100k-inhabitants figure (column D).
#
# Previously the 7-day window was CENTERED on the row's date (i.e. it
# summed 3 days before .. 3 days after). This edit re-centers the window
# so it ends on the current ("ultimo", last) day: the window becomes the
# 7 days up to and including the row's own date (r-6 .. r). Rows whose
# full trailing window doesn't fit inside the data range are left blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 184
$colNewCases = 2   # B: nuovi pos.
$colRolling = 3    # C: somma mobile 7gg.
$colPer100k = 4    # D: somma mobile 7gg. per 100mila abitanti
$population = 2141
$windowSize = 7

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $lo = $r - ($windowSize - 1)
    $hi = $r

    if ($lo -lt $firstRow -or $hi -gt $lastRow) {
        # Window doesn't fully fit in the data range -> leave blank,
        # but don't touch cells that are already blank.
        $existing = $ws.Cells.Item($r, $colRolling).Value2
        if ($null -ne $existing -and -not $existing.Equals("")) {
            $ws.Cells.Item($r, $colRolling).ClearContents()
            $ws.Cells.Item($r, $colPer100k).ClearContents()
        }
    } else {
        $sum = 0
        for ($i = $lo; $i -le $hi; $i++) {
            $sum += $ws.Cells.Item($i, $colNewCases).Value2
        }
        $ws.Cells.Item($r, $colRolling).Value = $sum
        $ws.Cells.Item($r, $colPer100k).Value = $sum * 100000 / $population
    }
}
